# Updated symbol list on Mon Jan  2 15:47:51 UTC 2023 with GitHub Actions
#
# The "Price" (column D) and "Volume(1h)" (column E) columns in this
# crypto-ranking sheet are stored as plain text (General/text cells),
# e.g. "246.78" or "0.95%". We force the NumberFormat to "@" (Text)
# before writing so Excel keeps the exact literal string instead of
# silently reinterpreting numeric- or percent-looking text as a real
# number/percentage value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $value) {
    $rng = $ws.Range($address)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2 - BNB
Set-TextValue "D2" "246.79"
Set-TextValue "E2" "0.99%"

# Row 3 - OKB
Set-TextValue "D3" "29.44"
Set-TextValue "E3" "7.23%"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.196"
Set-TextValue "E4" "2.82%"

# Row 5 - Cronos
Set-TextValue "D5" "0.05712"
Set-TextValue "E5" "0.53%"

# Row 6 - KuCoinToken
Set-TextValue "D6" "6.579"
Set-TextValue "E6" "1.68%"

# Row 7 - MXToken
Set-TextValue "D7" "0.8586"
Set-TextValue "E7" "4.58%"

# Row 8 - FTXToken
Set-TextValue "D8" "0.8810"
Set-TextValue "E8" "5.08%"

# Row 9 - WazirX
Set-TextValue "D9" "0.1369"
Set-TextValue "E9" "3.30%"

# Row 10 - MandalaExchangeToken
Set-TextValue "D10" "0.07070"
Set-TextValue "E10" "2.13%"

# Row 11 - BitrueCoin
Set-TextValue "D11" "0.02875"
Set-TextValue "E11" "0.39%"

# Row 12 - BitMartToken
Set-TextValue "D12" "0.09385"
Set-TextValue "E12" "-0.18%"

# Row 13 - BitForexToken
Set-TextValue "D13" "0.001521"
Set-TextValue "E13" "-0.03%"

# Row 14 - CoinExToken
Set-TextValue "D14" "0.04155"
Set-TextValue "E14" "1.10%"

# Row 15 - One
Set-TextValue "D15" "0.0005984"
Set-TextValue "E15" "-0.26%"

# Row 16 - TigerCash
Set-TextValue "D16" "0.006142"
Set-TextValue "E16" "0.12%"

# Row 17 - UpBots (price unchanged)
Set-TextValue "E17" "5,108.13%"

# Row 18 - LEO (price unchanged)
Set-TextValue "E18" "-0.80%"

# Row 19 - GateToken
Set-TextValue "D19" "3.065"
Set-TextValue "E19" "2.14%"

# Row 20 - BTSEToken
Set-TextValue "D20" "2.189"
Set-TextValue "E20" "-5.13%"

# Row 22 - LiechtensteinCryptoassetsExchange
Set-TextValue "D22" "0.03298"
Set-TextValue "E22" "3.24%"

# Row 23 - ProBitToken (price unchanged)
Set-TextValue "E23" "3.55%"

# Row 24 - MCDex
Set-TextValue "D24" "3.467"
Set-TextValue "E24" "-2.50%"

# Row 25 - ZBToken (price unchanged)
Set-TextValue "E25" "0.34%"

# Row 26 - HotbitToken
Set-TextValue "D26" "0.005056"
Set-TextValue "E26" "30.62%"

# Row 27 - BitKan
Set-TextValue "D27" "0.001220"
Set-TextValue "E27" "0.13%"

# Row 28 - NitroEx (price unchanged)
Set-TextValue "E28" "23.34%"

# Row 40 - IDEX
Set-TextValue "D40" "0.03752"
Set-TextValue "E40" "0.53%"

# Row 41 - KickToken
Set-TextValue "D41" "0.005653"
Set-TextValue "E41" "-7.61%"

# Row 42 - BKEXToken
Set-TextValue "D42" "0.1074"
Set-TextValue "E42" "1.95%"

# Row 43 - CEJI
Set-TextValue "D43" "0.002537"
Set-TextValue "E43" "10.32%"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.009972"
Set-TextValue "E44" "2.91%"

# Row 45 - CoinLion
Set-TextValue "D45" "0.00005128"
Set-TextValue "E45" "-1.69%"

# Row 47 - CoinbaseStockToken (volume unchanged)
Set-TextValue "D47" "0.07093"

# Row 48 - BOLO
Set-TextValue "D48" "0.002587"
Set-TextValue "E48" "0.65%"
